# Weekly update: insert the newest week's two price rows at the top of the
# data block (row 66 onward), pushing the existing history down by 2 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 66:145 down to 68:147, inserting two blank rows.
$ws.Rows("66:67").Insert()

# New row 66 - "Primera" quality, newest date.
$ws.Range("A66").Value = 2
$ws.Range("B66").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44629
$ws.Range("E66").Value = 4
$ws.Range("F66").Value = 100112043
$ws.Range("G66").Value = "Pepino ensalada"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 300
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 16000
$ws.Range("M66").Value = 15500
$ws.Range("N66").Value = "$/caja 70 unidades"
$ws.Range("O66").Value = "Provincia de Limarí"
$ws.Range("P66").Value = 221
$ws.Range("Q66").Value = 70
$ws.Range("R66").Value = "Hortaliza"

# New row 67 - "Segunda" quality, same newest date.
$ws.Range("A67").Value = 2
$ws.Range("B67").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C67").Value = "Coquimbo"
$ws.Range("D67").Value = 44629
$ws.Range("E67").Value = 4
$ws.Range("F67").Value = 100112043
$ws.Range("G67").Value = "Pepino ensalada"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Segunda"
$ws.Range("J67").Value = 200
$ws.Range("K67").Value = 13000
$ws.Range("L67").Value = 14000
$ws.Range("M67").Value = 13500
$ws.Range("N67").Value = "$/caja 100 unidades"
$ws.Range("O67").Value = "Provincia de Limarí"
$ws.Range("P67").Value = 135
$ws.Range("Q67").Value = 100
$ws.Range("R67").Value = "Hortaliza"
